$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Result" column header and "Pass" results for each data row
$ws.Range("D1").Value = "Result"
$ws.Range("D2").Value = "Pass"
$ws.Range("D3").Value = "Pass"
$ws.Range("D4").Value = "Pass"
$ws.Range("D5").Value = "Pass"

# Highlight the new results with a fill color (indexed color 11 / green)
$ws.Range("D2").Interior.Color = 65280
$ws.Range("D3").Interior.Color = 65280
$ws.Range("D4").Interior.Color = 65280
$ws.Range("D5").Interior.Color = 65280

# Match the selection left behind after the edit
[void]$ws.Range("D2:D5").Select()
